$d = $word.ActiveDocument

# 1. "ConverterApp" bullet: merge runs, text unchanged, but removes proofErr wrapping.
$d.Content.Find.Execute(
    "ConverterApp – a executable command line tool that generates meshlet data using DirectXMesh",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ConverterApp – a executable command line tool that generates meshlet data using DirectXMesh", 2)

# 2. "The ConverterApp project ..." paragraph: update FBX/OBJ wording and merge DirectXMesh/meshlet runs.
$d.Content.Find.Execute(
    "The ConverterApp project is a command line tool which can be used to generate meshlet data from an FBX or OBJ file. The tool leverages DirectXMesh integration of meshlet generation to generate meshlets from vertex & primitive data read from the input FBX file.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The ConverterApp project is a command line tool which can be used to generate meshlet data from an FBX file, OBJ file, or SDKMesh file. The tool leverages DirectXMesh integration of meshlet generation to generate meshlets from vertex & primitive data read from the input FBX file.", 2)

# 3. "-fz" bullet: merge, text unchanged.
$d.Content.Find.Execute(
    "-fz – Flips the Z axis of scene geometry. Default is false",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "-fz – Flips the Z axis of scene geometry. Default is false", 2)

# 4. "-i" bullet: merge, text unchanged.
$d.Content.Find.Execute(
    "-i – Forces vertex indices to 32-bits, even if 16-bits would suffice. Default is false",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "-i – Forces vertex indices to 32-bits, even if 16-bits would suffice. Default is false", 2)

# 5. "-t" bullet: merge Triangulates sentence, text unchanged.
$d.Content.Find.Execute(
    "Triangulates scene meshes file using the FbxGeometryConverter functionality",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Triangulates scene meshes file using the FbxGeometryConverter functionality", 2)

# 6. Update history: merge "DirectXMesh-like interface." into the 4/11/2020 entry.
$d.Content.Find.Execute(
    "4/11/2020 – Replaced meshlet generation interface with a thinner DirectXMesh-like interface.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "4/11/2020 – Replaced meshlet generation interface with a thinner DirectXMesh-like interface.", 2)

# Add a new "Update history" paragraph for the 10/17/2022 entry.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "10/17/2022 – Added support for reading from an SDKMesh file."

# 7. Footer copyright year 2021 -> 2022 (both footers).
$d.Content.Find.Execute("2021", $true, $false, $false, $false, $false, $true, 1, $false, "2022", 2)

foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Footers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute("2021", $true, $false, $false, $false, $false, $true, 1, $false, "2022", 2)
        }
    }
}
